$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1939393939393939
$ws.Range("C2").Value = 0.5393939393939394
$ws.Range("J2").Value = 0.02121212121212121
$ws.Range("P2").Value = 0.1212121212121212
$ws.Range("S2").Value = 0.1242424242424242
$ws.Range("B3").Value = 0.005291005291005291
$ws.Range("C3").Value = 0.02645502645502645
$ws.Range("J3").Value = 0.03703703703703703
$ws.Range("P3").Value = 0.7883597883597884
$ws.Range("S3").Value = 0.1428571428571428
$ws.Range("J4").Value = 0.103448275862069
$ws.Range("O4").Value = 0.01724137931034483
$ws.Range("P4").Value = 0.603448275862069
$ws.Range("S4").Value = 0.2758620689655172
$ws.Range("B6").Value = 0.05472636815920398
$ws.Range("D6").Value = 0.009950248756218905
$ws.Range("E6").Value = 0.009950248756218905
$ws.Range("F6").Value = 0.0845771144278607
$ws.Range("J6").Value = 0.3134328358208955
$ws.Range("O6").Value = 0.01990049751243781
$ws.Range("Q6").Value = 0.208955223880597
$ws.Range("R6").Value = 0.04477611940298507
$ws.Range("S6").Value = 0.2537313432835821
$ws.Range("B7").Value = 0.1415525114155251
$ws.Range("D7").Value = 0.0136986301369863
$ws.Range("F7").Value = 0.0547945205479452
$ws.Range("J7").Value = 0.1735159817351598
$ws.Range("O7").Value = 0.0365296803652968
$ws.Range("Q7").Value = 0.1735159817351598
$ws.Range("R7").Value = 0.0776255707762557
$ws.Range("S7").Value = 0.3287671232876712
$ws.Range("B8").Value = 0.08282828282828283
$ws.Range("D8").Value = 0.0202020202020202
$ws.Range("E8").Value = 0.00202020202020202
$ws.Range("F8").Value = 0.05656565656565657
$ws.Range("J8").Value = 0.1575757575757576
$ws.Range("O8").Value = 0.04242424242424243
$ws.Range("Q8").Value = 0.1878787878787879
$ws.Range("R8").Value = 0.08888888888888889
$ws.Range("S8").Value = 0.3616161616161616
$ws.Range("B9").Value = 0.08391608391608392
$ws.Range("D9").Value = 0.006993006993006993
$ws.Range("F9").Value = 0.0979020979020979
$ws.Range("J9").Value = 0.1118881118881119
$ws.Range("O9").Value = 0.02797202797202797
$ws.Range("Q9").Value = 0.2237762237762238
$ws.Range("R9").Value = 0.06993006993006994
$ws.Range("S9").Value = 0.3776223776223776
$ws.Range("B10").Value = 0.1230188679245283
$ws.Range("D10").Value = 0.03169811320754717
$ws.Range("F10").Value = 0.04603773584905661
$ws.Range("J10").Value = 0.1388679245283019
$ws.Range("O10").Value = 0.02792452830188679
$ws.Range("Q10").Value = 0.2332075471698113
$ws.Range("R10").Value = 0.06415094339622641
$ws.Range("S10").Value = 0.3350943396226415
$ws.Range("G11").Value = 0.1194029850746269
$ws.Range("J11").Value = 0.1104477611940299
$ws.Range("K11").Value = 0.1701492537313433
$ws.Range("L11").Value = 0.6
$ws.Range("G12").Value = 0.735
$ws.Range("J12").Value = 0.22
$ws.Range("K12").Value = 0.015
$ws.Range("L12").Value = 0.02
$ws.Range("S12").Value = 0.01
$ws.Range("G13").Value = 0.6956521739130435
$ws.Range("J13").Value = 0.2608695652173913
$ws.Range("S13").Value = 0.04347826086956522
$ws.Range("G14").Value = 0.75
$ws.Range("J14").Value = 0.25
$ws.Range("F15").Value = 0.03017241379310345
$ws.Range("H15").Value = 0.1163793103448276
$ws.Range("I15").Value = 0.05603448275862069
$ws.Range("J15").Value = 0.3405172413793103
$ws.Range("K15").Value = 0.07327586206896551
$ws.Range("M15").Value = 0.01293103448275862
$ws.Range("N15").Value = 0.008620689655172414
$ws.Range("O15").Value = 0.03448275862068965
$ws.Range("S15").Value = 0.3275862068965517
$ws.Range("F16").Value = 0.03225806451612903
$ws.Range("H16").Value = 0.1382488479262673
$ws.Range("I16").Value = 0.07834101382488479
$ws.Range("J16").Value = 0.4193548387096774
$ws.Range("K16").Value = 0.1382488479262673
$ws.Range("M16").Value = 0.01382488479262673
$ws.Range("O16").Value = 0.02764976958525346
$ws.Range("S16").Value = 0.152073732718894
$ws.Range("F17").Value = 0.02549019607843137
$ws.Range("H17").Value = 0.2019607843137255
$ws.Range("I17").Value = 0.07254901960784314
$ws.Range("J17").Value = 0.4019607843137255
$ws.Range("K17").Value = 0.09607843137254903
$ws.Range("M17").Value = 0.02745098039215686
$ws.Range("O17").Value = 0.0607843137254902
$ws.Range("S17").Value = 0.1137254901960784
$ws.Range("F18").Value = 0.01204819277108434
$ws.Range("H18").Value = 0.1686746987951807
$ws.Range("I18").Value = 0.04819277108433735
$ws.Range("J18").Value = 0.4457831325301205
$ws.Range("K18").Value = 0.09036144578313253
$ws.Range("M18").Value = 0.01204819277108434
$ws.Range("N18").Value = 0.006024096385542169
$ws.Range("O18").Value = 0.05421686746987952
$ws.Range("S18").Value = 0.1626506024096386
$ws.Range("F19").Value = 0.01144726083401472
$ws.Range("H19").Value = 0.2461161079313164
$ws.Range("I19").Value = 0.05396565821749796
$ws.Range("J19").Value = 0.3442354865085854
$ws.Range("K19").Value = 0.1316434995911692
$ws.Range("M19").Value = 0.02044153720359771
$ws.Range("N19").Value = 0.0008176614881439084
$ws.Range("O19").Value = 0.06541291905151267
$ws.Range("S19").Value = 0.1259198691741619
